# Updated cryptos list on Thu Nov 14 07:28:33 UTC 2024 with GitHub Actions
#
# Applies per-row Price (D) / Volume(1h) (E) refresh, plus two row swaps
# (rows 15/16 and rows 30/31) exactly as described by the target diff.
#
# Note: several new Price values look like plain decimal numbers
# (e.g. "216.14", "0.999", "1.00"). Excel's COM layer auto-converts such
# strings to numeric cell values when assigned directly, which would lose
# the original text formatting used throughout this sheet. Prefixing the
# value with a leading apostrophe forces Excel to store it as text (the
# apostrophe itself is not persisted in the cell's string content), just
# like manually typing e.g. '216.14 into a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.658.91"
$ws.Range("E2").Value = "  +2.87%  "

$ws.Range("D3").Value = "3.184.97"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'216.14"
$ws.Range("E5").Value = "  +5.99%  "

$ws.Range("D6").Value = "'622.84"
$ws.Range("E6").Value = "  +2.91%  "

$ws.Range("D7").Value = "'0.390"
$ws.Range("E7").Value = "  +4.48%  "

$ws.Range("D8").Value = "'0.689"
$ws.Range("E8").Value = "  +4.47%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "3.182.76"
$ws.Range("E10").Value = "  +1.32%  "

$ws.Range("D11").Value = "'0.569"
$ws.Range("E11").Value = "  +7.15%  "

$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("E13").Value = "  +5.94%  "

$ws.Range("E14").Value = "  +2.93%  "

# Row 15 / Row 16 swap (Avalanche <-> WrappedliquidstakedEther2.0), plus
# their own Price/Volume refresh.
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.773.46"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'33.10"
$ws.Range("E16").Value = "  +3.95%  "

$ws.Range("D17").Value = "89.357.17"
$ws.Range("E17").Value = "  +3.00%  "

$ws.Range("D18").Value = "3.200.80"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").Value = "'3.41"
$ws.Range("E19").Value = "  +14.27%  "

$ws.Range("E20").Value = "  +73.75%  "

$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").Value = "'435.38"
$ws.Range("E22").Value = "  +5.85%  "

$ws.Range("D23").Value = "'8.61"
$ws.Range("E23").Value = "  +2.05%  "

$ws.Range("D24").Value = "'5.07"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").Value = "'5.29"
$ws.Range("E25").Value = "  +3.17%  "

$ws.Range("D26").Value = "'11.93"
$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").Value = "'81.48"
$ws.Range("E27").Value = "  +11.44%  "

$ws.Range("D28").Value = "3.364.03"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("D29").Value = "'1.00"

# Row 30 / Row 31 swap (Cronos <-> Binance-PegBSC-USD), plus their own
# Price/Volume refresh.
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.158"
$ws.Range("E31").Value = "  -2.32%  "

$ws.Range("E32").Value = "  +35.71%  "

$ws.Range("D33").Value = "'544.04"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("E34").Value = "  +2.43%  "

$ws.Range("D35").Value = "'7.04"
$ws.Range("E35").Value = "  +6.90%  "

$ws.Range("E36").Value = "  +3.56%  "

$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("D38").Value = "'22.30"
$ws.Range("E38").Value = "  +2.43%  "

$ws.Range("D39").Value = "'22.37"
$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("D40").Value = "'0.128"
$ws.Range("E40").Value = "  -3.48%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("E43").Value = "  +2.20%  "

$ws.Range("D44").Value = "'0.372"
$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("D45").Value = "'146.09"
$ws.Range("E45").Value = "  -1.97%  "

$ws.Range("D46").Value = "'173.17"
$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("D47").Value = "'43.69"
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("D48").Value = "'0.758"
$ws.Range("E48").Value = "  +9.74%  "

$ws.Range("E49").Value = "  -2.06%  "

$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").Value = "'0.619"
$ws.Range("E51").Value = "  +6.02%  "
